$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.85
$ws.Activate()
$ws.Range("R10").Select() | Out-Null
